# The deck originally has slide 10 = "Hyper link" (sldId 263) and
# slide 11 = "Fin" (sldId 264). This edit removes the "Hyper link"
# slide entirely, so the former slide 11 ("Fin") becomes the new,
# final slide 10 and the deck goes from 11 slides down to 10.
$p = $ppt.ActivePresentation
$p.Slides.Item(10).Delete()
